$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update C12/D12: remove "Yes" mandatory flag, add description text
$ws.Range("C12").Value = ""
$ws.Range("D12").Value = "Only number is allowed, use ``null`` (without quotes) for empty value."

# Update C15/D15: remove "Yes" mandatory flag, add description text
$ws.Range("C15").Value = ""
$ws.Range("D15").Value = "Only number is allowed, use ``null`` (without quotes) for empty value."

# Update selection to D15
$ws.Range("D15").Select()
